$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 289
$ws1.Range("F4").Value = 1106
$ws1.Range("F5").Value = 578

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 289
$ws4.Range("F4").Value = 1106
$ws4.Range("F6").Value = 578
